$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "99÷7=14, 1"
$t.Cell(1,2).Range.Text = "17÷7=2, 3"
$t.Cell(1,3).Range.Text = "67÷9=7, 4"
$t.Cell(1,4).Range.Text = "71÷8=8, 7"
$t.Cell(1,5).Range.Text = "32÷3=10, 2"
$t.Cell(5,1).Range.Text = "10÷5=2, 0"
$t.Cell(5,2).Range.Text = "69÷9=7, 6"
$t.Cell(5,3).Range.Text = "54÷5=10, 4"
$t.Cell(5,4).Range.Text = "99÷3=33, 0"
$t.Cell(5,5).Range.Text = "95÷9=10, 5"
$t.Cell(9,1).Range.Text = "26÷4=6, 2"
$t.Cell(9,2).Range.Text = "58÷8=7, 2"
$t.Cell(9,3).Range.Text = "63÷5=12, 3"
$t.Cell(9,4).Range.Text = "81÷5=16, 1"
$t.Cell(9,5).Range.Text = "57÷9=6, 3"
$t.Cell(13,1).Range.Text = "87÷5=17, 2"
$t.Cell(13,2).Range.Text = "64÷4=16, 0"
$t.Cell(13,3).Range.Text = "19÷9=2, 1"
$t.Cell(13,4).Range.Text = "25÷2=12, 1"
$t.Cell(13,5).Range.Text = "19÷9=2, 1"
$t.Cell(17,1).Range.Text = "76÷4=19, 0"
$t.Cell(17,2).Range.Text = "94÷7=13, 3"
$t.Cell(17,3).Range.Text = "22÷4=5, 2"
$t.Cell(17,4).Range.Text = "18÷5=3, 3"
$t.Cell(17,5).Range.Text = "13÷2=6, 1"
